# Applies the "add scripts and results of eaglerc ngi mode" edit:
# - Renames the old HISAT-based method labels to their HISAT2 equivalents
# - Adds a new "Lasy-Seq (HISAT2)" column-group header
# - Adds "replicate" row labels above the two per-replicate tables

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1 header block: left block becomes "Lasy-Seq (HISAT2)", right block becomes
# "paired-end RNA-Seq (HISAT2)" (both previously referenced HISAT-only method names).
$ws.Range("C1").Value = "Lasy-Seq (HISAT2)"
$ws.Range("H1").Value = "paired-end RNA-Seq (HISAT2)"

# First data table's section label (row 3) also moves from HISAT to HISAT2.
$ws.Range("A3").Value = "paired-end RNA-Seq (HISAT2)"

# New "replicate" column labels above the #1..#6 rows in both tables.
$ws.Range("B2").Value = "replicate"

# B12/G12 were previously blank, centre-aligned placeholder cells; the new
# label text is written in the default (unformatted) style, so clear the old
# centring before writing the value.
$ws.Range("B12").ClearFormats()
$ws.Range("B12").Value = "replicate"
$ws.Range("G12").ClearFormats()
$ws.Range("G12").Value = "replicate"

# Leave the selection where the author left off editing.
$null = $ws.Range("I8").Select()
